$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Notes sheet: reword the "specific issue" note
# ---------------------------------------------------------------------------
$notes = $wb.Worksheets.Item("Notes")
$notes.Range("A3").Value = "Specific issue: survey_IDs are not unique within a study"

# ---------------------------------------------------------------------------
# studies sheet: rename header, rename existing row, add a second study row
# ---------------------------------------------------------------------------
$studies = $wb.Worksheets.Item("studies")
$studies.Range("A1").Value = "study_id"
$studies.Range("A2").Value = "study01"
$studies.Range("B2").Value = "example name"

$studiesRow3 = @("study02", "example name", "other", "Blaggs_etal", 2024, "https://doi.org/10.1093%2Fgenetics%2F16.2.97")
$col = 1
foreach ($v in $studiesRow3) {
    $studies.Cells.Item(3, $col).Value = $v
    $col = $col + 1
}
$studies.Hyperlinks.Add($studies.Range("F3"), "https://doi.org/10.1093%2Fgenetics%2F16.2.97")
$studies.Range("F3").Style = "Hyperlink"

# ---------------------------------------------------------------------------
# surveys sheet: rename headers, recolor header font black, add 3 more rows
# ---------------------------------------------------------------------------
$surveys = $wb.Worksheets.Item("surveys")
$surveys.Range("B1").Value = "survey_id"
$surveys.Range("E1").Value = "latitude"
$surveys.Range("F1").Value = "longitude"
$surveys.Range("A2").Value = "study01"

$surveys.Range("A1:G1").Font.Color = 0
$surveys.Range("K1").Font.Color = 0
$surveys.Range("H1:J1").Font.Color = 0

$surveysRows = @(
    @("study01", "S02", "Gambia", "example site", 0, 0, "example data", "2020-01-01", "2020-01-01", "2020-01-01", "example data"),
    @("study02", "S01", "Mali",   "example site", 0, 0, "example data", "2020-01-01", "2020-01-01", "2020-01-01", "example data"),
    @("study02", "S01", "Mali",   "example site", 0, 0, "example data", "2020-01-01", "2020-01-01", "2020-01-01", "example data")
)
$r = 3
foreach ($row in $surveysRows) {
    $col = 1
    foreach ($v in $row) {
        $cell = $surveys.Cells.Item($r, $col)
        $cell.Value = $v
        if ($col -ge 8 -and $col -le 10) {
            $cell.Style = "Normal"
            $cell.NumberFormat = "@"
        }
        $col = $col + 1
    }
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# counts sheet has no content changes (only shared-string reindexing, which
# follows automatically from the edits above)
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# View state: active tab moves from "surveys" to "studies"
# ---------------------------------------------------------------------------
$studies.Activate()
$studies.Range("A2").Select()
$surveys.Range("D12").Select()
